$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSEval")

$ws.Range("A1").Value = "Input [Raw Data File]"
$ws.Range("K1").Value = "Output [Derived Data File]"
